$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Data: week label (Dn), start date (En), end date (Fn) for rows 5..14
$weeks = @(
    @{ Row = 5;  Label = "S 1";  Start = "01/01/2021"; End = "01/07/2021" },
    @{ Row = 6;  Label = "S 2";  Start = "01/08/2021"; End = "01/14/2021" },
    @{ Row = 7;  Label = "S 3";  Start = "01/15/2021"; End = "01/21/2021" },
    @{ Row = 8;  Label = "S 4";  Start = "01/22/2021"; End = "01/28/2021" },
    @{ Row = 9;  Label = "S 5";  Start = "01/29/2021"; End = "02/04/2021" },
    @{ Row = 10; Label = "S 6";  Start = "02/05/2021"; End = "02/11/2021" },
    @{ Row = 11; Label = "S 7";  Start = "02/12/2021"; End = "02/18/2021" },
    @{ Row = 12; Label = "S 8";  Start = "02/19/2021"; End = "02/25/2021" },
    @{ Row = 13; Label = "S 9";  Start = "02/26/2021"; End = "03/04/2021" },
    @{ Row = 14; Label = "S 10"; Start = "03/05/2021"; End = "03/11/2021" }
)

foreach ($w in $weeks) {
    $r = $w.Row
    $ws.Range("D$r").Value = $w.Label
    $ws.Range("E$r").Value = $w.Start
    $ws.Range("F$r").Value = $w.End
}

# Column F: widen and drop the auto bestFit sizing
$ws.Columns.Item(6).ColumnWidth = 11.7

# Update the active selection shown when the file is opened
$ws.Range("F17").Select()
